$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Good Morning" (R10's greeting, cell E8) was replaced with "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the author's last selection (E8) in the saved view state
[void]$ws.Range("E8").Select()
